$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A352").Value = 44807
$ws.Range("B352").Value = "HP63C1093"
$ws.Range("C352").Value = "SX4"
$ws.Range("D352").Value = "PMS"
$ws.Range("E352").Value = "WORK DONE DELIVERED"
$ws.Range("F352").Value = 7335
$ws.Range("G352").Value = "GPAY"

$ws.Range("A353").Value = 44807
$ws.Range("B353").Value = "KA25P8050"
$ws.Range("C353").Value = "NANO"
$ws.Range("D353").Value = "CLUTCH PROBLEM"
$ws.Range("E353").Value = "WORK DONE DELIVERED"
$ws.Range("F353").Value = 6914

$ws.Range("A354").Value = 44807
$ws.Range("B354").Value = "TN14P4378"
$ws.Range("C354").Value = "CELERIO"
$ws.Range("D354").Value = "GENERAL CHECKUP "
$ws.Range("E354").Value = "WORK DONE DELIVERED"
$ws.Range("F354").Value = 1500
$ws.Range("G354").Value = "P PAY"

$ws.Range("A355").Value = 44807
$ws.Range("B355").Value = "KL43B1476"
$ws.Range("C355").Value = "ALTO"
$ws.Range("D355").Value = "PMS"
$ws.Range("E355").Value = "WORK IN PROGRESS"

$ws.Range("A356").Value = 44776
$ws.Range("B356").Value = "KA03MU8528"
$ws.Range("C356").Value = "SWIFT"
$ws.Range("D356").Value = "BODY SHOP"
$ws.Range("E356").Value = "WORK IN PROGRESS"

$ws.Range("A357").Value = 44807
$ws.Range("B357").Value = "KA53MA1192"
$ws.Range("C357").Value = "ETIOS LIVA"
$ws.Range("D357").Value = "SUSPENSION"
$ws.Range("E357").Value = "WORK IN PROGRESS"

$ws.Range("A358").Value = 44777
$ws.Range("B358").Value = "JH09F7221"
$ws.Range("C358").Value = "SWIFT"
$ws.Range("D358").Value = "CLUTCH PROBLEM"
$ws.Range("E358").Value = "WORK IN PROGRESS"

$ws.Range("A359").Value = 44808
$ws.Range("B359").Value = "KA10M2591"
$ws.Range("C359").Value = "POLO"
$ws.Range("D359").Value = "CLUTCH PROBLEM  & PMS"
$ws.Range("E359").Value = "WORK IN PROGRESS"

$ws.Range("A360").Value = 44809
$ws.Range("B360").Value = "KA02AF2947"
$ws.Range("C360").Value = "SWIFT DZIRE"
$ws.Range("D360").Value = "STARTING PROBLEM"
$ws.Range("E360").Value = "WORK DONE DELIVERED"
$ws.Range("F360").Value = 3800
$ws.Range("G360").Value = "CASH"

$ws.Range("A361").Value = 44809
$ws.Range("B361").Value = "KA11B6109"
$ws.Range("C361").Value = "XYLO"
$ws.Range("D361").Value = "STARTING PROBLEM"
$ws.Range("E361").Value = "WORK IN PROGRESS"

$ws.Range("A362").Value = 44809
$ws.Range("B362").Value = "KL49E1440"
$ws.Range("C362").Value = "FIESTA"
$ws.Range("D362").Value = "GEAR STEERING CHANGE"
$ws.Range("E362").Value = "WORK IN PROGRESS"

$ws.Range("A363").Value = 44809
$ws.Range("B363").Value = "UP14CK8538"
$ws.Range("C363").Value = "SAEL"
$ws.Range("D363").Value = "BODY SHOP"
$ws.Range("E363").Value = "WORK IN PROGRESS"

$ws.Range("A364").Value = 44809
$ws.Range("B364").Value = "KA22P5434"
$ws.Range("C364").Value = "PUNTO"
$ws.Range("D364").Value = "WIPER BLADE CHANGE"
$ws.Range("E364").Value = "WORK DONE DELIVERED"
$ws.Range("F364").Value = 820
$ws.Range("G364").Value = "P PAY"

$ws.Range("A365").Value = 44809
$ws.Range("B365").Value = "KA03ML5436"
$ws.Range("C365").Value = "H-CITY"
$ws.Range("D365").Value = "GENERAL CHECKUP"
$ws.Range("E365").Value = "WORK DONE DELIVERED"
$ws.Range("F365").Value = 5113
$ws.Range("G365").Value = "CREDIT"

$ws.Range("A366").Value = 44809
$ws.Range("B366").Value = "KA03NA2866"
$ws.Range("C366").Value = "POLO"
$ws.Range("D366").Value = "PMS"
$ws.Range("E366").Value = "WORK DONE DELIVERED"
$ws.Range("F366").Value = 10244

$ws.Range("A367").Value = 44809
$ws.Range("B367").Value = "KA53N2941"
$ws.Range("C367").Value = "CAMRY"
$ws.Range("D367").Value = "AC REFLLING "
$ws.Range("E367").Value = "WORK DONE DELIVERED"
$ws.Range("F367").Value = 2796
$ws.Range("G367").Value = "GPAY"

$ws.Range("A368").Value = 44809
$ws.Range("B368").Value = "KA03MZ2202"
$ws.Range("C368").Value = "AUDI A4"
$ws.Range("D368").Value = "AC REFLLING "
$ws.Range("E368").Value = "WORK DONE DELIVERED"
$ws.Range("F368").Value = 2950
$ws.Range("G368").Value = "GPAY"

$ws.Range("H368").Select()
